# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E3) and
# "Correspond Handback DateTime" (H3) for the d41011ff... row on
# both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 14:50:58"
$wsZhCn.Range("H3").Value = "2016-03-22 14:51:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 14:51:02"
$wsDeDe.Range("H3").Value = "2016-03-22 14:51:36"
